$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 5 ("BPMN") and Slide 6 ("ACTIVITY DIAGRAM") swap their title text and
# their picture's crop / position / size / name / alt-text, while staying in
# the same position within the deck.
# ---------------------------------------------------------------------------

$s5 = $p.Slides.Item(5)
$s6 = $p.Slides.Item(6)

# --- Title text swap ---
$s5.Shapes.Item(3).TextFrame.TextRange.Text = "ACTIVITY DIAGRAM"
$s6.Shapes.Item(3).TextFrame.TextRange.Text = "BPMN"

# --- Picture on slide 5: take on slide 6's picture's current geometry/crop ---
$pic5 = $s5.Shapes.Item(10)
$pic5.Name = "Picture 4"
$pic5.AlternativeText = "Diagram`n`nDescription automatically generated"
$pic5.PictureFormat.CropLeft = 527.877105
$pic5.Left = 26.888739585876465
$pic5.Top = 67.34621810913086
$pic5.Width = 906.2224426269532
$pic5.Height = 413.0130767822266

# --- Picture on slide 6: take on slide 5's (original) picture's geometry/crop ---
$pic6 = $s6.Shapes.Item(10)
$pic6.Name = "Picture 2"
$pic6.AlternativeText = "Diagram, schematic`n`nDescription automatically generated"
$pic6.PictureFormat.CropLeft = 543.508035
$pic6.Left = 29.916142463684086
$pic6.Top = 69.29023361206055
$pic6.Width = 900.1676940917969
$pic6.Height = 412.82685852050787

# ---------------------------------------------------------------------------
# Slide 7 ("TAHAPAN PROTOTYPE"): retitle and nudge the two header divider
# connectors.
# ---------------------------------------------------------------------------

$s7 = $p.Slides.Item(7)

$conn7 = $s7.Shapes.Item(2)   # "Straight Connector 7" (right-hand divider)
$conn7.Left = 665.6470642089844
$conn7.Width = 294.3529205322266

$s7.Shapes.Item(3).TextFrame.TextRange.Text = "TAHAPAN SD:C PROTOTYPE"

$conn13 = $s7.Shapes.Item(4)  # "Straight Connector 13" (left-hand divider)
$conn13.Left = -22.875117301940914
$conn13.Top = 42.110864639282234
